$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 "标签" (tag): shorten "传奇巫师" (Legendary Wizard) to just "巫师" (Wizard)
$ws.Range("D2").Value = "巫师"

# I2 "引言" (quote): replace the short quote with the longer, repeated one
# (note: the literal two characters "\n" are part of the text itself, not an
# actual line break)
$ws.Range("I2").Value = "牢大，我想你了\n牢大牢大，我想你了牢大，我想你了牢大，我想你了牢大，我想你了牢大，我想你了牢大，我想你了牢大，我想你了"

# The author's selection ended on the quote cell they just edited
$ws.Range("I2").Select()
